# 00 - 04 fig 1 updates
# Appends 9 new rows (28-36) of figure metadata to the sheet, matching the
# commit that extended the used range from A1:I27 to A1:I36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=28; A="Burgelman_2002_ASQ_Strategy as Vector and the Inertia of Coevolutionary Lock-in.pdf"; B="F1_P28_Burgelman_2002_ASQ_Strategy as Vector and the Inertia of Coevolutionary Lock-in.png"; C=1; D=2002; E=29; F=0; G="final_figures/2002/F1_P28_Burgelman_2002_ASQ_Strategy as Vector and the Inertia of Coevolutionary Lock-in.png"; H="conceptual diagram"; I=$true },
    @{ Row=29; A="Dutton et al._2002_OrgSci_Red Ligh Green Light.pdf"; B="F1_P4_Dutton et al._2002_OrgSci_Red Ligh Green Light.png"; C=1; D=2002; E=5; F=0; G="final_figures/2002/F1_P4_Dutton et al._2002_OrgSci_Red Ligh Green Light.png"; H="none selected"; I=$true },
    @{ Row=30; A="Fitzgerald, Ferlie, Wood & Hawkins_2002_HR_Interlocking interactions, the diffusion of innovations in health care.pdf"; B="F1_P6_Fitzgerald, Ferlie, Wood & Hawkins_2002_HR_Interlocking interactions, the diffusion of innovations in health care.png"; C=1; D=2002; E=7; F=0; G="final_figures/2002/F1_P6_Fitzgerald, Ferlie, Wood & Hawkins_2002_HR_Interlocking interactions, the diffusion of innovations in health care.png"; H="2x2 matrix"; I=$true },
    @{ Row=31; A="Foreman & Whetten_2002_OrgSci_Members identification with Multiple Identity Orgs_Quant.pdf"; B="F1_P3_Foreman & Whetten_2002_OrgSci_Members identification with Multiple Identity Orgs_Quant.png"; C=1; D=2002; E=4; F=0; G="final_figures/2002/F1_P3_Foreman & Whetten_2002_OrgSci_Members identification with Multiple Identity Orgs_Quant.png"; H="conceptual diagram"; I=$true },
    @{ Row=32; A="Hodgkinson & Wright_2002_OrgSci_Confronting Strategic Inertia in a Top Management Team - Learning from Failure.pdf"; B="F1_P14_Hodgkinson & Wright_2002_OrgSci_Confronting Strategic Inertia in a Top Management Team - Learning from Failure.png"; C=1; D=2002; E=15; F=0; G="final_figures/2002/F1_P14_Hodgkinson & Wright_2002_OrgSci_Confronting Strategic Inertia in a Top Management Team - Learning from Failure.png"; H="process diagram"; I=$true },
    @{ Row=33; A="Kitchener_2002_OrgStudies_Mobilitzing the Logic of Managerialism in Professional Fields.pdf"; B="F1_P9_Kitchener_2002_OrgStudies_Mobilitzing the Logic of Managerialism in Professional Fields.png"; C=1; D=2002; E=10; F=0; G="final_figures/2002/F1_P9_Kitchener_2002_OrgStudies_Mobilitzing the Logic of Managerialism in Professional Fields.png"; H="process diagram"; I=$true },
    @{ Row=34; A="Repenning & Sterman_2002_ASQ_Capability Traps and Self-Confirming Attribution Errors.pdf"; B="F1_P11_Repenning & Sterman_2002_ASQ_Capability Traps and Self-Confirming Attribution Errors.png"; C=1; D=2002; E=12; F=0; G="final_figures/2002/F1_P11_Repenning & Sterman_2002_ASQ_Capability Traps and Self-Confirming Attribution Errors.png"; H="process diagram"; I=$true },
    @{ Row=35; A="Snell & Tseng_2002_OrgSci_Moral Atmosphere and Moral Influence under China's Network Capitalism.pdf"; B="F1_P2_Snell & Tseng_2002_OrgSci_Moral Atmosphere and Moral Influence under China's Network Capitalism.png"; C=1; D=2002; E=3; F=-90; G="final_figures/2002/F1_P2_Snell & Tseng_2002_OrgSci_Moral Atmosphere and Moral Influence under China's Network Capitalism.png"; H="conceptual diagram"; I=$true },
    @{ Row=36; A="Snell_2002_OrgSci_The Learning Organization, SEnsegiving and Psychological Contracts - a Hong Kong Case.pdf"; B="F1_P15_Snell_2002_OrgSci_The Learning Organization, SEnsegiving and Psychological Contracts - a Hong Kong Case.png"; C=1; D=2002; E=16; F=-90; G="final_figures/2002/F1_P15_Snell_2002_OrgSci_The Learning Organization, SEnsegiving and Psychological Contracts - a Hong Kong Case.png"; H="conceptual diagram"; I=$true }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
}
